# Apply updated classification metrics (D/E columns) and flipped
# Success booleans (C column) for rows 8-10 and 18-20, plus the
# recomputed Cross Entropy Loss / Success % pairs (F/G) on rows 11 and 21.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = 0.537691799692685
$ws.Range("E2").Value = 0.537691799692685

# Row 3
$ws.Range("D3").Value = 0.008298101627418562
$ws.Range("E3").Value = 0.008298101627418562

# Row 4
$ws.Range("D4").Value = 0.00001052450147654965
$ws.Range("E4").Value = 0.00001052450147654965

# Row 5
$ws.Range("D5").Value = 0.07155001002993018
$ws.Range("E5").Value = 0.07155001002993018

# Row 6
$ws.Range("D6").Value = 0.926531823764956
$ws.Range("E6").Value = 0.926531823764956

# Row 7
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.107023546700823
$ws.Range("E8").Value = 0.8929764532991771

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.04280406929257809
$ws.Range("E9").Value = 0.9571959307074219

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.000005433649152542796
$ws.Range("E10").Value = 0.9999945663508475

# Row 11
$ws.Range("D11").Value = 0.9999901441969345
$ws.Range("E11").Value = 0.00000985580306545053
$ws.Range("F11").Value = 2.097374677658081
$ws.Range("G11").Value = 0.5

# Row 12
$ws.Range("D12").Value = 0.6818550405123067
$ws.Range("E12").Value = 0.6818550405123067

# Row 13
$ws.Range("D13").Value = 0.2296225764641644
$ws.Range("E13").Value = 0.2296225764641644

# Row 14
$ws.Range("D14").Value = 0.000000002222882800240327
$ws.Range("E14").Value = 0.000000002222882800240327

# Row 15
$ws.Range("D15").Value = 0.002509176330194134
$ws.Range("E15").Value = 0.002509176330194134

# Row 16
$ws.Range("D16").Value = 0.9425429279696622
$ws.Range("E16").Value = 0.9425429279696622

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.0000000301669996619711
$ws.Range("E18").Value = 0.9999999698330003

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.02019003278735198
$ws.Range("E19").Value = 0.979809967212648

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.000000002742611859183358
$ws.Range("E20").Value = 0.9999999972573882

# Row 21
$ws.Range("D21").Value = 0.9999999732107685
$ws.Range("E21").Value = 0.00000002678923149979084
$ws.Range("F21").Value = 4.51987886428833
$ws.Range("G21").Value = 0.5

